# auxUnits.xlsx demo fix ("incomplete fix of demo")
#
# The "combustEff" variable-name label was renamed to the more readable
# "combustion eff" everywhere it is used as a header label:
#   - sheet5 "c Power Station" cell F3
#   - sheet6 "v Power Station" cell C1
#
# Once both occurrences of the old shared string are gone, it is dropped
# from the workbook's shared-string table automatically and every other
# (unrelated) shared-string index that sorted after it shifts down by one
# - that ripple is purely mechanical and falls out of the two value edits
# below, it does not need to be reproduced by hand.

$wb = $excel.ActiveWorkbook

$wsPowerC      = $wb.Worksheets.Item("c Power Station")
$wsPowerV      = $wb.Worksheets.Item("v Power Station")
$wsCaptureVals = $wb.Worksheets.Item("v CO2 Capture")

# Rename the label in both places it occurs.
$wsPowerC.Range("F3").Value = "combustion eff"
$wsPowerV.Range("C1").Value = "combustion eff"

# Update the active selection on the sheet that doesn't change (matches the
# recorded selection state in the saved file).
$wsCaptureVals.Range("E16").Select() | Out-Null

# Move the live selection to reflect where the edit was made, and leave the
# workbook open on the "c Power Station" tab (it becomes the active sheet).
$wsPowerV.Range("C1").Select() | Out-Null
$wsPowerC.Range("F3").Select() | Out-Null
